# Apply "Add 2022-06-24 data" updates to the
# fonds-solidarite-volet-1-regional-classe-effectif dataset.
#
# For each listed row, update column C (nombre_aides) and column E
# (montant_total) to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 65;  C = 2020;   E = 13691790 },
    @{ Row = 91;  C = 151201; E = 482924639 },
    @{ Row = 92;  C = 409263; E = 1597003409 },
    @{ Row = 93;  C = 209648; E = 1309917695 },
    @{ Row = 94;  C = 94226;  E = 918566141 },
    @{ Row = 96;  C = 17317;  E = 796835901 },
    @{ Row = 98;  C = 812;    E = 117970793 },
    @{ Row = 107; C = 6392;   E = 21970305 },
    @{ Row = 144; C = 24419;  E = 202157062 },
    @{ Row = 153; C = 99155;  E = 337748796 },
    @{ Row = 158; C = 3848;   E = 140904137 },
    @{ Row = 174; C = 226108; E = 900747512 },
    @{ Row = 175; C = 80786;  E = 486199250 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
